$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the formula in B27 to add the extra time entries
$ws.Range("B27").Formula = "= 4.5 + 4 + 1.5 + 1 + 2.5 + 1 + 1.5 + 2.5 + 3.25 + 1.5"

# Recalculate the workbook so dependent formulas (E2, E3, G3, H3) update
$excel.Calculate()

# Update the selected cell/range as recorded in the sheet view
$ws.Range("F23").Select()
